$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 23:22"

# Estados Unidos (row 4) - new totals
$ws.Range("B4").Value = 462391
$ws.Range("C4").Value = 27464
$ws.Range("D4").Value = 25139
$ws.Range("E4").Value = 420798
$ws.Range("G4").Value = 1666
$ws.Range("H4").Value = 16454

# Alemania (row 8) - new totals
$ws.Range("B8").Value = 116801
$ws.Range("C8").Value = 3505
$ws.Range("D8").Value = 52407
$ws.Range("E8").Value = 61943

# Row 43 used to hold Serbia's figures; the country list now places
# "Emiratos Arabes Unidos" ahead of "Serbia", so row 43 becomes the
# (freshly updated) Emiratos Arabes Unidos entry ...
$ws.Range("A43").Value = "Emiratos Arabes Unidos"
$ws.Range("B43").Value = 2990
$ws.Range("C43").Value = 331
$ws.Range("D43").Value = 268
$ws.Range("E43").Value = 2708
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 14

# ... and row 44 becomes Serbia, carrying forward its previous figures.
$ws.Range("A44").Value = "Serbia"
$ws.Range("B44").Value = 2867
$ws.Range("C44").Value = 201
$ws.Range("D44").Value = 118
$ws.Range("E44").Value = 2683
$ws.Range("F44").Value = 127
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 66

# Colombia (row 50) - new totals
$ws.Range("B50").Value = 2223
$ws.Range("C50").Value = 169
$ws.Range("D50").Value = 174
$ws.Range("E50").Value = 1980
$ws.Range("G50").Value = 14
$ws.Range("H50").Value = 69
